# Confidence regions for experiments - add "Experiment 3" sheet, mirroring the
# layout of "Experiment 1" / "Experiment 2", and update the selection left on
# "Experiment 2" so it no longer points at the old working cells.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Experiment 1")
$ws2 = $wb.Worksheets.Item("Experiment 2")

# ---------------------------------------------------------------------------
# 1. Experiment 2: move the leftover selection off the scratch column (K) and
#    onto the main table, and drop the "this is the active tab" marker (the
#    new sheet added below becomes the active one instead).
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A11:G29").Select()

# ---------------------------------------------------------------------------
# 2. Add the new "Experiment 3" worksheet after "Experiment 2".
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Experiment 3"

# ---------------------------------------------------------------------------
# 3. Bring over the formatting used on "Experiment 1" / "Experiment 2" so the
#    new sheet's number formats / fonts match exactly (copy format only).
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

$ws1.Range("A1:A9").Copy()
$ws3.Range("A1:A9").PasteSpecial($xlPasteFormats)

$ws1.Range("B12:B21").Copy()
$ws3.Range("B12:B21").PasteSpecial($xlPasteFormats)

$ws1.Range("B22:B29").Copy()
$ws3.Range("B22:B29").PasteSpecial($xlPasteFormats)

$ws1.Range("C12:C29").Copy()
$ws3.Range("C12:C29").PasteSpecial($xlPasteFormats)

# Leftover formatted-but-empty cells in columns N / S (mirrors the stray
# formatting left on column K of "Experiment 2" / columns I,J,M,N of
# "Experiment 1").
$ws1.Range("M16:M31").Copy()
$ws3.Range("N14:N29").PasteSpecial($xlPasteFormats)

$ws1.Range("N11:N26").Copy()
$ws3.Range("S14:S29").PasteSpecial($xlPasteFormats)

$ws1.Range("N27:N28").Copy()
$ws3.Range("S30:S31").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# 4. Raw data table (rows 1-9).
# ---------------------------------------------------------------------------
$ws3.Range("B1").Value = "workedm"
$ws3.Range("C1").Value = "morekids"
$ws3.Range("D1").Value = "girls2"
$ws3.Range("E1").Value = "freq"

$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = 0
$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 52474

$ws3.Range("A3").Value = 2
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 0
$ws3.Range("D3").Value = 1
$ws3.Range("E3").Value = 15030

$ws3.Range("A4").Value = 3
$ws3.Range("B4").Value = 0
$ws3.Range("C4").Value = 1
$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 38678

$ws3.Range("A5").Value = 4
$ws3.Range("B5").Value = 0
$ws3.Range("C5").Value = 1
$ws3.Range("D5").Value = 1
$ws3.Range("E5").Value = 13959

$ws3.Range("A6").Value = 5
$ws3.Range("B6").Value = 1
$ws3.Range("C6").Value = 0
$ws3.Range("D6").Value = 0
$ws3.Range("E6").Value = 70211

$ws3.Range("A7").Value = 6
$ws3.Range("B7").Value = 1
$ws3.Range("C7").Value = 0
$ws3.Range("D7").Value = 1
$ws3.Range("E7").Value = 20027

$ws3.Range("A8").Value = 7
$ws3.Range("B8").Value = 1
$ws3.Range("C8").Value = 1
$ws3.Range("D8").Value = 0
$ws3.Range("E8").Value = 32345

$ws3.Range("A9").Value = 8
$ws3.Range("B9").Value = 1
$ws3.Range("C9").Value = 1
$ws3.Range("D9").Value = 1
$ws3.Range("E9").Value = 11930

# ---------------------------------------------------------------------------
# 5. Z=0 block (rows 11-16).
# ---------------------------------------------------------------------------
$ws3.Range("A11").Value = "Z=0"
$ws3.Range("B11").Formula = "=SUM(E2,E4,E6,E8)"

$ws3.Range("A12").Value = "f(0,0)"
$ws3.Range("B12").Formula = "=E2"
$ws3.Range("C12").Formula = "=B12/`$B`$11"

$ws3.Range("A13").Value = "f(0,1)"
$ws3.Range("B13").Formula = "=E4"
$ws3.Range("C13").Formula = "=B13/`$B`$11"

$ws3.Range("A14").Value = "f(1,0)"
$ws3.Range("B14").Formula = "=E6"
$ws3.Range("C14").Formula = "=B14/`$B`$11"

$ws3.Range("A15").Value = "f(1,1)"
$ws3.Range("B15").Formula = "=E8"
$ws3.Range("C15").Formula = "=B15/`$B`$11"

# ---------------------------------------------------------------------------
# 6. Z=1 block (rows 17-22).
# ---------------------------------------------------------------------------
$ws3.Range("A17").Value = "Z=1"
$ws3.Range("B17").Formula = "=E3+E5+E7+E9"

$ws3.Range("A18").Value = "f(0,0)"
$ws3.Range("B18").Formula = "=E3"
$ws3.Range("C18").Formula = "=B18/`$B`$17"

$ws3.Range("A19").Value = "f(0,1)"
$ws3.Range("B19").Formula = "=E5"
$ws3.Range("C19").Formula = "=B19/`$B`$17"

$ws3.Range("A20").Value = "f(1,0)"
$ws3.Range("B20").Formula = "=E7"
$ws3.Range("C20").Formula = "=B20/`$B`$17"

$ws3.Range("A21").Value = "f(1,1)"
$ws3.Range("B21").Formula = "=E9"
$ws3.Range("C21").Formula = "=B21/`$B`$17"

# ---------------------------------------------------------------------------
# 7. ATE(neg) confidence block (rows 23-26).
# ---------------------------------------------------------------------------
$ws3.Range("A23").Value = "ATE(neg)"
$ws3.Range("B23").Value = "LB"
$ws3.Range("C23").Value = "UB"
$ws3.Range("E23").Value = "Confidence"
$ws3.Range("F23").Value = "LB"
$ws3.Range("G23").Value = "UB"

$ws3.Range("A24").Value = "1-Q(0)"
$ws3.Range("B24").Formula = "=MAX(C14+C15,C20+C21)"
$ws3.Range("C24").Formula = "=MIN(C13+C14+C15,C19+C20+C21)"
$ws3.Range("E24").Value = "1-Q(0)"
$ws3.Range("F24").Value = 0.52689520000000001
$ws3.Range("G24").Value = 0.73107279999999997

$ws3.Range("A25").Value = "1-Q(1)"
$ws3.Range("B25").Formula = "=MAX(C15,C21)"
$ws3.Range("C25").Formula = "=MIN(C15+C14,C21+C20)"
$ws3.Range("E25").Value = "1-Q(1)"
$ws3.Range("F25").Value = 0.1927372
$ws3.Range("G25").Value = 0.52887930000000005

# ---------------------------------------------------------------------------
# 8. ATE(pos) confidence block (rows 27-29).
# ---------------------------------------------------------------------------
$ws3.Range("A27").Value = "ATE(pos)"
$ws3.Range("B27").Value = "LB"
$ws3.Range("C27").Value = "UB"
$ws3.Range("E27").Value = "Confidence"
$ws3.Range("F27").Value = "LB"
$ws3.Range("G27").Value = "UB"

$ws3.Range("A28").Value = "1-Q(0)"
$ws3.Range("B28").Formula = "=MAX(C14,C20)"
$ws3.Range("C28").Formula = "=MIN(C14+C15,C20+C21)"
$ws3.Range("E28").Value = "1-Q(0)"
$ws3.Range("F28").Value = 0.3603287
$ws3.Range("G28").Value = 0.52887930000000005

$ws3.Range("A29").Value = "1-Q(1)"
$ws3.Range("B29").Formula = "=MAX(C15+C14,C21+C20)"
$ws3.Range("C29").Formula = "=MIN(C12+C14+C15,C18+C20+C21)"
$ws3.Range("E29").Value = "1-Q(1)"
$ws3.Range("F29").Value = 0.52689520000000001
$ws3.Range("G29").Value = 0.77417729999999996

# ---------------------------------------------------------------------------
# 9. Leave the same kind of leftover selection on the new (now active) sheet
#    that the other two sheets have.
# ---------------------------------------------------------------------------
$ws3.Activate()
$ws3.Range("K27").Select()
